$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A110:H110").EntireRow.Insert()

$ws.Cells.Item(110, 1).Value = "R11"
$ws.Cells.Item(110, 2).Value = "CZ"
$ws.Cells.Item(110, 3).Value = 1805
$ws.Cells.Item(110, 4).Value = 1805.00001305901
$ws.Cells.Item(110, 5).Value = 1804.9999993369
$ws.Cells.Item(110, 6).Value = 9081752.44566366
$ws.Cells.Item(110, 7).Value = 0.0000130590051412582
$ws.Cells.Item(110, 8).Value = 0.000000663101673126221
